$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date cell T2 (copy the date style from a neighboring date cell so it
# reuses the existing numFmtId 14 style instead of minting a new one), then set
# its value to the date serial for 2019-09-18.
$ws.Range("U2").Copy($ws.Range("T2"))
$ws.Range("T2").Value = 43726

# Convert W2 from a text date ("18-09-2019") to a real date serial value,
# keeping its existing date-formatted style.
$ws.Range("W2").Value = 43726

# Update the view: clear the frozen/scrolled topLeftCell and move the
# selection to L2.
[void]$ws.Range("L2").Select()
